$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Make sure Sheet1's own view is reset to the top-left (A1) before we
# switch the active tab away from it.
$ws1.Activate() | Out-Null
$ws1.Range("A1").Select() | Out-Null

# Insert the new "readme" sheet right after Sheet1.
$readme = $wb.Worksheets.Add($null, $ws1)
$readme.Name = "readme"

# Header row (bold).
$readme.Range("A1").Value = "variable"
$readme.Range("B1").Value = "description"
$readme.Range("A1:B1").Font.Bold = $true

# Column-description rows, pairing each original data-column header with
# a human readable description of what it contains.
$readme.Range("A2").Value = "stream"
$readme.Range("B2").Value = "stream network"

$readme.Range("A3").Value = "site"
$readme.Range("B3").Value = "site name"

$readme.Range("A4").Value = "latitude"
$readme.Range("B4").Value = "latitude (degrees)"

$readme.Range("A5").Value = "longitude"
$readme.Range("B5").Value = "longitude (degrees)"

$readme.Range("A6").Value = "datetime_0"
$readme.Range("B6").Value = "measurement start time"

$readme.Range("A7").Value = "datetime_1"
$readme.Range("B7").Value = "measurement end time"

$readme.Range("A8").Value = "mean_flux"
$readme.Range("B8").Value = "mean co2 flux in umol/m2/sec"

$readme.Range("A9").Value = "n"
$readme.Range("B9").Value = "number of observations"

$readme.Range("A10").Value = "ph"
$readme.Range("B10").Value = "ph"

$readme.Range("A11").Value = "temp"
$readme.Range("B11").Value = "degrees celcius"

$readme.Range("A12").Value = "alk"
$readme.Range("B12").Value = "meq/L"

$readme.Range("A13").Value = "index"
$readme.Range("B13").Value = "qpoint id"

# Make the new readme sheet the active tab/selection, matching the
# workbook's saved "activeTab" pointing at it.
$readme.Activate() | Out-Null
$readme.Range("B4").Select() | Out-Null
